# Update "Generate Report for Handback" timestamps.
# The handoff/handback datetimes recorded for the 5ee1c237... (.md) file
# are refreshed to reflect a regenerated handback report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 08:54:11"
$wsZhCn.Range("H2").Value = "2016-03-19 08:54:27"

# de-de sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 08:54:14"
$wsDeDe.Range("H2").Value = "2016-03-19 08:54:32"
